## EngineMerge.xlsx consolidation edit
## - Rename shared string "ATT" -> "AP_ATT" (DQ1 on Tabelle2)
## - Insert 6 new parameter columns (COM1, COM1_NAME, COM2, COM2_NAME, NAV1_F, NAV2_F)
##   right before the old "END_OF_COL" column on Tabelle2
## - Re-point the selection like the saved workbook did

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# 1) Rename ATT -> AP_ATT (this is the existing DQ1 header cell)
$ws.Range("DQ1").Value = "AP_ATT"

# 2) Insert 6 blank columns before the old DR column (END_OF_COL), shifting
#    END_OF_COL/Title from DR:DS to DX:DY. The new columns inherit the
#    formatting of the column to their left (DQ), matching style ids 11/12.
$ws.Range("DR1:DW40").EntireColumn.Insert()

# 3) New header row (row 1) values for the inserted columns
$ws.Cells.Item(1, 122).Value = "COM1"
$ws.Cells.Item(1, 123).Value = "COM1_NAME"
$ws.Cells.Item(1, 124).Value = "COM2"
$ws.Cells.Item(1, 125).Value = "COM2_NAME"
$ws.Cells.Item(1, 126).Value = "NAV1_F"
$ws.Cells.Item(1, 127).Value = "NAV2_F"

# 4) Fill the new columns on every aircraft data row (2-40) with the same
#    "|" (Any) placeholder used throughout the rest of the row.
for ($row = 2; $row -le 40; $row++) {
    for ($col = 122; $col -le 127; $col++) {
        $ws.Cells.Item($row, $col).Value = "|"
    }
}

# 5) Match the saved selection state
$ws.Range("BN42").Select() | Out-Null
